{"js": "// Replace the 25 two-digit multiplication problems in the practice table.\n// Each \"before\" value is unique in the document, so a targeted search+replace\n// for each pair is safe and will not cascade into already-replaced cells.\nconst replacements = [\n  [\"57\u00d740=\", \"70\u00d790=\"],\n  [\"76\u00d798=\", \"69\u00d782=\"],\n  [\"25\u00d798=\", \"98\u00d792=\"],\n  [\"26\u00d795=\", \"97\u00d791=\"],\n  [\"82\u00d755=\", \"55\u00d720=\"],\n  [\"40\u00d714=\", \"37\u00d790=\"],\n  [\"80\u00d754=\", \"27\u00d786=\"],\n  [\"15\u00d748=\", \"40\u00d795=\"],\n  [\"72\u00d779=\", \"59\u00d738=\"],\n  [\"78\u00d788=\", \"45\u00d798=\"],\n  [\"64\u00d784=\", \"34\u00d725=\"],\n  [\"80\u00d746=\", \"20\u00d773=\"],\n  [\"75\u00d755=\", \"12\u00d798=\"],\n  [\"78\u00d759=\", \"81\u00d727=\"],\n  [\"76\u00d743=\", \"99\u00d716=\"],\n  [\"95\u00d732=\", \"81\u00d735=\"],\n  [\"11\u00d763=\", \"66\u00d754=\"],\n  [\"78\u00d742=\", \"25\u00d750=\"],\n  [\"84\u00d747=\", \"82\u00d784=\"],\n  [\"76\u00d711=\", \"93\u00d728=\"],\n  [\"92\u00d744=\", \"70\u00d737=\"],\n  [\"13\u00d751=\", \"13\u00d770=\"],\n  [\"98\u00d768=\", \"75\u00d754=\"],\n  [\"99\u00d776=\", \"29\u00d782=\"],\n  [\"95\u00d790=\", \"25\u00d726=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Pattern not found: ${before}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 25 two-digit multiplication problems in the practice table.\n# Each \"before\" value is unique in the document, so a targeted\n# Find/Replace for each pair is safe and will not cascade into cells\n# that were already updated earlier in the loop.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Before = \"57\u00d740=\"; After = \"70\u00d790=\" },\n    @{ Before = \"76\u00d798=\"; After = \"69\u00d782=\" },\n    @{ Before = \"25\u00d798=\"; After = \"98\u00d792=\" },\n    @{ Before = \"26\u00d795=\"; After = \"97\u00d791=\" },\n    @{ Before = \"82\u00d755=\"; After = \"55\u00d720=\" },\n    @{ Before = \"40\u00d714=\"; After = \"37\u00d790=\" },\n    @{ Before = \"80\u00d754=\"; After = \"27\u00d786=\" },\n    @{ Before = \"15\u00d748=\"; After = \"40\u00d795=\" },\n    @{ Before = \"72\u00d779=\"; After = \"59\u00d738=\" },\n    @{ Before = \"78\u00d788=\"; After = \"45\u00d798=\" },\n    @{ Before = \"64\u00d784=\"; After = \"34\u00d725=\" },\n    @{ Before = \"80\u00d746=\"; After = \"20\u00d773=\" },\n    @{ Before = \"75\u00d755=\"; After = \"12\u00d798=\" },\n    @{ Before = \"78\u00d759=\"; After = \"81\u00d727=\" },\n    @{ Before = \"76\u00d743=\"; After = \"99\u00d716=\" },\n    @{ Before = \"95\u00d732=\"; After = \"81\u00d735=\" },\n    @{ Before = \"11\u00d763=\"; After = \"66\u00d754=\" },\n    @{ Before = \"78\u00d742=\"; After = \"25\u00d750=\" },\n    @{ Before = \"84\u00d747=\"; After = \"82\u00d784=\" },\n    @{ Before = \"76\u00d711=\"; After = \"93\u00d728=\" },\n    @{ Before = \"92\u00d744=\"; After = \"70\u00d737=\" },\n    @{ Before = \"13\u00d751=\"; After = \"13\u00d770=\" },\n    @{ Before = \"98\u00d768=\"; After = \"75\u00d754=\" },\n    @{ Before = \"99\u00d776=\"; After = \"29\u00d782=\" },\n    @{ Before = \"95\u00d790=\"; After = \"25\u00d726=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Before\n    $find.Replacement.Text = $pair.After\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    [void]$find.Execute($pair.Before, $false, $false, $false, $false, $false, $true, 1, $false, $pair.After, 2)\n}\n"}
